# Bugs.xlsx - "Finished testing Game class"
# Adds the results of testing the Game class: a fix description for the
# existing "Player health drops below zero" bug, plus a newly discovered
# bug (HP reaching 0 doesn't end the game) with its own fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix column for the previously-logged "Player health drops below
#     zero" bug (row 8) ---
$ws.Range("C8").Value = "HP is set to zero if incoming damage is greater than HP"

# --- New bug entry (row 9) ---
$ws.Range("A9").Value = "Game does not end if a player's HP reached 0"
$ws.Range("B9").Value = "Player is asked for turn input before HP check is done to end the game"
$ws.Range("C9").Value = "Moved HP check code to execute before player turn starts"

# Match the explicit row heights used by the rest of the sheet.
$ws.Rows(9).RowHeight = 15.75
$ws.Rows(10).RowHeight = 15.75

# --- Widen the Bug/Description columns to fit the new, longer text ---
$ws.Range("A1").ColumnWidth = 76.5885416666667
$ws.Range("B1").ColumnWidth = 99.5885416666667
